$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new Wins/Losses/Ties columns (AD, AE, AF) in row 1,
# matching the formatting already used by the other header cells in row 1
# (bold font, thin box border, centered horizontally, top-aligned vertically).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop

$headerRange.Borders.LineStyle = 1   # xlContinuous
$headerRange.Borders.Weight = 2      # xlThin

# Fill in the team record values for every data row (2 through 46).
$ws.Range("AD2:AD46").Value = 75
$ws.Range("AE2:AE46").Value = 87
$ws.Range("AF2:AF46").Value = 0
